$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = '44.212.60'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = '  +2.57%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = '2.426.18'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = '  +2.17%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = '  -0.09%  '
$ws.Range("E5").Value2 = '  +1.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '100.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = '  +3.99%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '0.514'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = '  +1.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '0.502'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = '  +0.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '35.28'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = '  +3.21%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '0.0800'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = '  +1.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '18.91'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = '  +3.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '0.123'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = '  +1.97%  '
$ws.Range("E14").Value2 = '  +1.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '2.803.82'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = '  +2.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '2.469.90'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = '  +2.52%  '
$ws.Range("E17").Value2 = '  +3.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '44.149.63'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = '  +2.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '12.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = '  +0.87%  '
$ws.Range("E20").Value2 = '  +1.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = '0.0₃0906'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = '  +2.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '68.60'
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value2 = '  +5.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '240.56'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = '  +2.21%  '
$ws.Range("E25").Value2 = '  +1.57%  '
$ws.Range("E26").Value2 = '  +0.05%  '
$ws.Range("E27").Value2 = '  +1.85%  '
$ws.Range("E28").Value2 = '  -1.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '9.58'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = '  +4.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '32.86'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = '  +4.95%  '
$ws.Range("E31").Value2 = '  +12.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '18.69'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = '  +7.45%  '
$ws.Range("E33").Value2 = '  +2.25%  '
$ws.Range("E34").Value2 = '  -0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = '0.0759'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = '  +1.11%  '
$ws.Range("B36").Value2 = 'Monero'
$ws.Range("C36").Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '131.63'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = '  +25.82%  '
$ws.Range("B37").Value2 = 'ARBITRUM'
$ws.Range("C37").Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '1.90'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = '  +3.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '4.48'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = '  +4.47%  '
$ws.Range("E39").Value2 = '  +3.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '2.29'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = '  -1.35%  '
$ws.Range("E41").Value2 = '  +1.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '21.61'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = '  -2.98%  '
$ws.Range("E43").Value2 = '  +2.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '1.951.24'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = '  -0.33%  '
$ws.Range("E45").Value2 = '  +2.00%  '
$ws.Range("E46").Value2 = '  +4.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '9.40'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = '  +2.48%  '
$ws.Range("E48").Value2 = '  +9.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '53.47'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = '  +1.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '73.76'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = '  +2.52%  '
$ws.Range("E51").Value2 = '  +1.09%  '
